$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Text edit: "Địa chỉ liên lạc khi cần báo tin: <<ThiSinh_DienThoai>>"
#    becomes  "Địa chỉ liên lạc khi cần báo tin: <<ThiSinh_DCNhanGiayBao>>"
#    (only the FIRST occurrence of the placeholder in the document - there
#    is a second, unrelated "ThiSinh_DienThoai" placeholder further down
#    that must stay untouched).
# ---------------------------------------------------------------------------

# Anchor on text that only occurs once, right before the run we need to edit.
$anchor = $d.Content
$anchor.Find.Execute("cần báo tin", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorEnd = $anchor.End

# Restrict the search to a small window right after the anchor and locate
# the ": <<" literal that precedes the placeholder field.
$searchRng = $d.Range($anchorEnd, $anchorEnd + 80)
$searchRng.Find.Execute(": <<", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$colonStart = $searchRng.Start
$colonEnd = $searchRng.End

# Split the merged ": <<" run into ": " and "<<" (matching the target,
# which keeps them as two separate runs with identical formatting).
# Toggling Bold on/off forces a run split without leaving any formatting
# residue behind.
$splitRng = $d.Range($colonStart + 2, $colonEnd)
$splitRng.Bold = $true
$splitRng.Bold = $false

# Replace the placeholder field name and colour it explicit black, matching
# the target run's added <w:color w:val="000000"/>.
$fieldRng = $d.Range($colonEnd, $colonEnd + 40)
$fieldRng.Find.Execute("ThiSinh_DienThoai", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fieldRng.Text = "ThiSinh_DCNhanGiayBao"
$fieldRng.Font.Color = 0

# ---------------------------------------------------------------------------
# 2) styles.xml latent-style table bookkeeping: add the "Table *" latent
#    style exceptions (w:uiPriority/semiHidden/unhideWhenUsed) that Word
#    writes once a style is touched/saved through a newer Word build.
# ---------------------------------------------------------------------------

$tableLatentStyles = @(
    "Normal Table",
    "Table Simple 1",
    "Table Simple 2",
    "Table Simple 3",
    "Table Classic 1",
    "Table Classic 2",
    "Table Classic 3",
    "Table Classic 4",
    "Table Colorful 1",
    "Table Colorful 2",
    "Table Colorful 3",
    "Table Columns 1",
    "Table Columns 2",
    "Table Columns 3",
    "Table Columns 4",
    "Table Columns 5",
    "Table Grid 1",
    "Table Grid 2",
    "Table Grid 3",
    "Table Grid 4",
    "Table Grid 5",
    "Table Grid 6",
    "Table Grid 7",
    "Table Grid 8",
    "Table List 1",
    "Table List 2",
    "Table List 3",
    "Table List 4",
    "Table List 5",
    "Table List 6",
    "Table List 7",
    "Table List 8",
    "Table 3D effects 1",
    "Table 3D effects 2",
    "Table 3D effects 3",
    "Table Contemporary",
    "Table Elegant",
    "Table Professional",
    "Table Subtle 1",
    "Table Subtle 2",
    "Table Web 1",
    "Table Web 2",
    "Table Web 3",
    "Table Theme"
)

foreach ($name in $tableLatentStyles) {
    $ls = $d.LatentStyles.Add($name)
    $ls.SemiHidden = $true
    $ls.UnhideWhenUsed = $true
}

# ---------------------------------------------------------------------------
# 3) Style definitions gain explicit uiPriority + unhideWhenUsed entries.
# ---------------------------------------------------------------------------

$dpf = $d.Styles("Default Paragraph Font")
$dpf.Priority = 1
$dpf.UnhideWhenUsed = $true

$tblNormal = $d.Styles("Normal Table")
$tblNormal.Priority = 99
$tblNormal.UnhideWhenUsed = $true

$noList = $d.Styles("No List")
$noList.Priority = 99
$noList.UnhideWhenUsed = $true

# ---------------------------------------------------------------------------
# 4) Fix the stray leading space in the custom "Char" style name.
# ---------------------------------------------------------------------------

$charStyle = $d.Styles(" Char")
$charStyle.NameLocal = "Char"

Write-Output "edit complete"
